$wb = $excel.ActiveWorkbook

# Worksheets: "dust_event" (first sheet) and "dust_event_small" (second sheet)
$wsDust = $wb.Worksheets.Item("dust_event")
$wsSmall = $wb.Worksheets.Item("dust_event_small")

# Update the value in dust_event!E8: was "IDL Small Batch Processor failed",
# now should read "IDL Batch Processor failed" (a new shared string).
$wsDust.Range("E8").Value = "IDL Batch Processor failed"

# Update the selection on dust_event_small to A9, scrolled back to A1.
$wsSmall.Activate()
$wsSmall.Range("A1").Select()
$wsSmall.Range("A9").Select()

# Make dust_event the active/selected sheet with selection E18.
$wsDust.Activate()
$wsDust.Range("E18").Select()
